# Update Case_3_22 vm_pu results: slack bus voltage setpoint changed from 1.05 to 1.02 pu
# (recomputed load-flow results for all buses, rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049304338986681
$ws.Range("D2").Value = 1.055130053900059
$ws.Range("E2").Value = 1.056501405637298
$ws.Range("F2").Value = 1.067202617919128
$ws.Range("I2").Value = 1.041436550823106
$ws.Range("J2").Value = 1.05434349526658
$ws.Range("K2").Value = 1.057871359157953
$ws.Range("L2").Value = 1.059238943010632
$ws.Range("M2").Value = 1.069911112862725
$ws.Range("N2").Value = 1.021861909914845

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.050486103456782
$ws.Range("D3").Value = 1.056045993168506
$ws.Range("E3").Value = 1.057525543812355
$ws.Range("F3").Value = 1.068247721612781
$ws.Range("I3").Value = 1.041675432624431
$ws.Range("J3").Value = 1.055173348696137
$ws.Range("K3").Value = 1.058600416955825
$ws.Range("L3").Value = 1.060076196937312
$ws.Range("M3").Value = 1.070771385219403
$ws.Range("N3").Value = 1.0221448054848

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051250912047018
$ws.Range("D4").Value = 1.05663858498241
$ws.Range("E4").Value = 1.058188613512241
$ws.Range("F4").Value = 1.068924267406847
$ws.Range("I4").Value = 1.041828644991222
$ws.Range("J4").Value = 1.055709907325331
$ws.Range("K4").Value = 1.059071453325637
$ws.Range("L4").Value = 1.060617730851739
$ws.Range("M4").Value = 1.07132772551549
$ws.Range("N4").Value = 1.022327544860663

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051572469230515
$ws.Range("D5").Value = 1.05688769089877
$ws.Range("E5").Value = 1.058467460053979
$ws.Range("F5").Value = 1.06920875768336
$ws.Range("I5").Value = 1.041892729815014
$ws.Range("J5").Value = 1.055935378247961
$ws.Range("K5").Value = 1.059269306483603
$ws.Range("L5").Value = 1.060845337957471
$ws.Range("M5").Value = 1.07156153605152
$ws.Range("N5").Value = 1.022404293565194

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.051626461982366
$ws.Range("D6").Value = 1.056929515723559
$ws.Range("E6").Value = 1.05851428501706
$ws.Range("F6").Value = 1.069256528984793
$ws.Range("I6").Value = 1.041903470841286
$ws.Range("J6").Value = 1.055973230051373
$ws.Range("K6").Value = 1.059302516902507
$ws.Range("L6").Value = 1.060883551045848
$ws.Range("M6").Value = 1.071600789479079
$ws.Range("N6").Value = 1.022417175616448

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051255208584186
$ws.Range("D7").Value = 1.05664191362522
$ws.Range("E7").Value = 1.058192339111534
$ws.Range("F7").Value = 1.068928068503022
$ws.Range("I7").Value = 1.041829502575421
$ws.Range("J7").Value = 1.05571292046479
$ws.Range("K7").Value = 1.059074097719159
$ws.Range("L7").Value = 1.060620772359985
$ws.Range("M7").Value = 1.071330849997578
$ws.Range("N7").Value = 1.022328570675032

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.049703695940944
$ws.Range("D8").Value = 1.05543961686185
$ws.Range("E8").Value = 1.056847438548839
$ws.Range("F8").Value = 1.067555754818856
$ws.Range("I8").Value = 1.041517563648592
$ws.Range("J8").Value = 1.054624033622296
$ws.Range("K8").Value = 1.058117895304649
$ws.Range("L8").Value = 1.059521943678522
$ws.Range("M8").Value = 1.0702019109033
$ws.Range("N8").Value = 1.02195758059959

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046970666558614
$ws.Range("D9").Value = 1.053320387438705
$ws.Range("E9").Value = 1.054480484128981
$ws.Range("F9").Value = 1.065139808757819
$ws.Range("I9").Value = 1.040957471813429
$ws.Range("J9").Value = 1.052702105667583
$ws.Range("K9").Value = 1.056427480145685
$ws.Range("L9").Value = 1.057583931728902
$ws.Range("M9").Value = 1.068210169731933
$ws.Range("N9").Value = 1.021301451550812

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.04514920837732
$ws.Range("D10").Value = 1.05190714216134
$ws.Range("E10").Value = 1.052904472651981
$ws.Range("F10").Value = 1.063530678909091
$ws.Range("I10").Value = 1.040577075020172
$ws.Range("J10").Value = 1.051418666046747
$ws.Range("K10").Value = 1.055296852807829
$ws.Range("L10").Value = 1.056290741214123
$ws.Range("M10").Value = 1.06688071377593
$ws.Range("N10").Value = 1.020862419226139

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044360614401918
$ws.Range("D11").Value = 1.05129508880144
$ws.Range("E11").Value = 1.052222502828815
$ws.Range("F11").Value = 1.06283426085895
$ws.Range("I11").Value = 1.040410697873046
$ws.Range("J11").Value = 1.050862405671896
$ws.Range("K11").Value = 1.054806401884194
$ws.Range("L11").Value = 1.055730490639314
$ws.Range("M11").Value = 1.066304655149787
$ws.Range("N11").Value = 1.020671930185967

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.044067710480903
$ws.Range("D12").Value = 1.05106772833402
$ws.Range("E12").Value = 1.051969256571661
$ws.Range("F12").Value = 1.062575631764826
$ws.Range("I12").Value = 1.040348648081341
$ws.Range("J12").Value = 1.050655706689905
$ws.Range("K12").Value = 1.054624093727138
$ws.Range("L12").Value = 1.055522344509273
$ws.Range("M12").Value = 1.066090621532854
$ws.Range("N12").Value = 1.020601116115069

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.044130538710998
$ws.Range("D13").Value = 1.051116498681163
$ws.Range("E13").Value = 1.052023575697011
$ws.Range("F13").Value = 1.062631106255473
$ws.Range("I13").Value = 1.04036196928276
$ws.Range("J13").Value = 1.050700047912858
$ws.Range("K13").Value = 1.05466320546311
$ws.Range("L13").Value = 1.055566994564877
$ws.Range("M13").Value = 1.066136535182267
$ws.Range("N13").Value = 1.020616308598018

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.04433640257494
$ws.Range("D14").Value = 1.051276295454308
$ws.Range("E14").Value = 1.052201568035692
$ws.Range("F14").Value = 1.062812881447347
$ws.Range("I14").Value = 1.040405573916545
$ws.Range("J14").Value = 1.050845321486906
$ws.Range("K14").Value = 1.054791334947192
$ws.Range("L14").Value = 1.055713286118056
$ws.Range("M14").Value = 1.066286964286262
$ws.Range("N14").Value = 1.020666077855698

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04446324403939
$ws.Range("D15").Value = 1.051374749313387
$ws.Range("E15").Value = 1.052311243947743
$ws.Range("F15").Value = 1.062924885965747
$ws.Range("I15").Value = 1.040432407052179
$ws.Range("J15").Value = 1.050934818829424
$ws.Range("K15").Value = 1.054870262136814
$ws.Range("L15").Value = 1.055803415311623
$ws.Range("M15").Value = 1.066379640687149
$ws.Range("N15").Value = 1.02069673465482

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045201546944372
$ws.Range("D16").Value = 1.0519477598131
$ws.Range("E16").Value = 1.052949742261274
$ws.Range("F16").Value = 1.063576905135066
$ws.Range("I16").Value = 1.040588081881912
$ws.Range("J16").Value = 1.051455572154751
$ws.Range("K16").Value = 1.055329383810409
$ws.Range("L16").Value = 1.056327917038031
$ws.Range("M16").Value = 1.06691893652567
$ws.Range("N16").Value = 1.02087505323361

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045664693128134
$ws.Range("D17").Value = 1.052307164850368
$ws.Range("E17").Value = 1.053350376501545
$ws.Range("F17").Value = 1.063985991849323
$ws.Range("I17").Value = 1.040685287340159
$ws.Range("J17").Value = 1.051782086630205
$ws.Range("K17").Value = 1.05561714236143
$ws.Range("L17").Value = 1.056656844994723
$ws.Range("M17").Value = 1.067257116418874
$ws.Range("N17").Value = 1.020986804584963

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045934849117461
$ws.Range("D18").Value = 1.052516789244284
$ws.Range("E18").Value = 1.053584103465989
$ws.Range("F18").Value = 1.064224638597728
$ws.Range("I18").Value = 1.040741825206882
$ws.Range("J18").Value = 1.051972486621644
$ws.Range("K18").Value = 1.055784901996755
$ws.Range("L18").Value = 1.056848675104729
$ws.Range("M18").Value = 1.067454332902305
$ws.Range("N18").Value = 1.021051950102973

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.046026967179896
$ws.Range("D19").Value = 1.052588263974715
$ws.Range("E19").Value = 1.053663805761779
$ws.Range("F19").Value = 1.064306016646477
$ws.Range("I19").Value = 1.040761075948995
$ws.Range("J19").Value = 1.052037399571748
$ws.Range("K19").Value = 1.055842089280632
$ws.Range("L19").Value = 1.056914079479538
$ws.Range("M19").Value = 1.067521572171399
$ws.Range("N19").Value = 1.02107415674195

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045615000833812
$ws.Range("D20").Value = 1.052268605177269
$ws.Range("E20").Value = 1.053307387712724
$ws.Range("F20").Value = 1.063942097265389
$ws.Range("I20").Value = 1.040674874709314
$ws.Range("J20").Value = 1.05174705993213
$ws.Range("K20").Value = 1.055586277392764
$ws.Range("L20").Value = 1.056621557037557
$ws.Range("M20").Value = 1.067220836863111
$ws.Range("N20").Value = 1.02097481856628

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044275780395757
$ws.Range("D21").Value = 1.051229239747627
$ws.Range("E21").Value = 1.052149151862987
$ws.Range("F21").Value = 1.062759351763178
$ws.Range("I21").Value = 1.040392740338213
$ws.Range("J21").Value = 1.050802544228909
$ws.Range("K21").Value = 1.054753607677806
$ws.Range("L21").Value = 1.055670208123609
$ws.Range("M21").Value = 1.0662426683182
$ws.Range("N21").Value = 1.020651423650527

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043433845385305
$ws.Range("D22").Value = 1.050575653295622
$ws.Range("E22").Value = 1.05142131448577
$ws.Range("F22").Value = 1.062016010949686
$ws.Range("I22").Value = 1.040213905306916
$ws.Range("J22").Value = 1.050208231705815
$ws.Range("K22").Value = 1.054229306389496
$ws.Range("L22").Value = 1.055071801796567
$ws.Range("M22").Value = 1.065627308788162
$ws.Range("N22").Value = 1.020447757045545

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043880163392736
$ws.Range("D23").Value = 1.050922141037642
$ws.Range("E23").Value = 1.051807117891026
$ws.Range("F23").Value = 1.062410041836186
$ws.Range("I23").Value = 1.040308846230601
$ws.Range("J23").Value = 1.050523331635309
$ws.Range("K23").Value = 1.054507321416255
$ws.Range("L23").Value = 1.055389052688934
$ws.Range("M23").Value = 1.065953555457612
$ws.Range("N23").Value = 1.02055575638935

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045637454609283
$ws.Range("D24").Value = 1.052286028665967
$ws.Range("E24").Value = 1.053326812360435
$ws.Range("F24").Value = 1.063961931235065
$ws.Range("I24").Value = 1.040679580224614
$ws.Range("J24").Value = 1.051762887144699
$ws.Range("K24").Value = 1.05560022420616
$ws.Range("L24").Value = 1.056637502232753
$ws.Range("M24").Value = 1.067237230150322
$ws.Range("N24").Value = 1.020980234646927

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.04767711633378
$ws.Range("D25").Value = 1.053868333852516
$ws.Range("E25").Value = 1.0550920524656
$ws.Range("F25").Value = 1.065764123505895
$ws.Range("I25").Value = 1.041103503215741
$ws.Range("J25").Value = 1.053199348148521
$ws.Range("K25").Value = 1.056865141558537
$ws.Range("L25").Value = 1.058085161545365
$ws.Range("M25").Value = 1.068725368813876
$ws.Range("N25").Value = 1.021471360835934
